$wb = $excel.ActiveWorkbook

# --- Sheet "Fixture frances": randomized fixture re-generated by new run.py / minimax pipeline ---
$wsFixture = $wb.Worksheets.Item("Fixture frances")
$wsFixture.Range("B2").Value = "@BOL"
$wsFixture.Range("C2").Value = "ARG"
$wsFixture.Range("D2").Value = "VEN"
$wsFixture.Range("E2").Value = "@COL"
$wsFixture.Range("F2").Value = "PAR"
$wsFixture.Range("G2").Value = "@PER"
$wsFixture.Range("H2").Value = "@ECU"
$wsFixture.Range("I2").Value = "CHI"
$wsFixture.Range("J2").Value = "URU"
$wsFixture.Range("K2").Value = "@ARG"
$wsFixture.Range("L2").Value = "@VEN"
$wsFixture.Range("M2").Value = "COL"
$wsFixture.Range("N2").Value = "@PAR"
$wsFixture.Range("O2").Value = "PER"
$wsFixture.Range("P2").Value = "ECU"
$wsFixture.Range("Q2").Value = "@CHI"
$wsFixture.Range("R2").Value = "@URU"
$wsFixture.Range("S2").Value = "BOL"
$wsFixture.Range("C3").Value = "@BRA"
$wsFixture.Range("D3").Value = "@URU"
$wsFixture.Range("E3").Value = "CHI"
$wsFixture.Range("H3").Value = "COL"
$wsFixture.Range("K3").Value = "BRA"
$wsFixture.Range("L3").Value = "URU"
$wsFixture.Range("M3").Value = "@CHI"
$wsFixture.Range("P3").Value = "@COL"
$wsFixture.Range("B4").Value = "VEN"
$wsFixture.Range("C4").Value = "@PER"
$wsFixture.Range("D4").Value = "@PAR"
$wsFixture.Range("E4").Value = "BRA"
$wsFixture.Range("F4").Value = "@URU"
$wsFixture.Range("G4").Value = "CHI"
$wsFixture.Range("H4").Value = "@ARG"
$wsFixture.Range("J4").Value = "@ECU"
$wsFixture.Range("K4").Value = "PER"
$wsFixture.Range("L4").Value = "PAR"
$wsFixture.Range("M4").Value = "@BRA"
$wsFixture.Range("N4").Value = "URU"
$wsFixture.Range("O4").Value = "@CHI"
$wsFixture.Range("P4").Value = "ARG"
$wsFixture.Range("R4").Value = "ECU"
$wsFixture.Range("S4").Value = "@VEN"
$wsFixture.Range("B5").Value = "CHI"
$wsFixture.Range("C5").Value = "@VEN"
$wsFixture.Range("D5").Value = "ARG"
$wsFixture.Range("E5").Value = "@ECU"
$wsFixture.Range("F5").Value = "COL"
$wsFixture.Range("G5").Value = "@BOL"
$wsFixture.Range("H5").Value = "PAR"
$wsFixture.Range("J5").Value = "@BRA"
$wsFixture.Range("K5").Value = "VEN"
$wsFixture.Range("L5").Value = "@ARG"
$wsFixture.Range("M5").Value = "ECU"
$wsFixture.Range("N5").Value = "@COL"
$wsFixture.Range("O5").Value = "BOL"
$wsFixture.Range("P5").Value = "@PAR"
$wsFixture.Range("R5").Value = "BRA"
$wsFixture.Range("S5").Value = "@CHI"
$wsFixture.Range("B6").Value = "@URU"
$wsFixture.Range("D6").Value = "ECU"
$wsFixture.Range("E6").Value = "@ARG"
$wsFixture.Range("F6").Value = "VEN"
$wsFixture.Range("G6").Value = "@COL"
$wsFixture.Range("H6").Value = "PER"
$wsFixture.Range("I6").Value = "@BRA"
$wsFixture.Range("J6").Value = "BOL"
$wsFixture.Range("L6").Value = "@ECU"
$wsFixture.Range("M6").Value = "ARG"
$wsFixture.Range("N6").Value = "@VEN"
$wsFixture.Range("O6").Value = "COL"
$wsFixture.Range("P6").Value = "@PER"
$wsFixture.Range("Q6").Value = "BRA"
$wsFixture.Range("R6").Value = "@BOL"
$wsFixture.Range("S6").Value = "URU"
$wsFixture.Range("B7").Value = "@PAR"
$wsFixture.Range("C7").Value = "COL"
$wsFixture.Range("D7").Value = "BOL"
$wsFixture.Range("E7").Value = "@VEN"
$wsFixture.Range("F7").Value = "@ECU"
$wsFixture.Range("G7").Value = "BRA"
$wsFixture.Range("H7").Value = "@CHI"
$wsFixture.Range("K7").Value = "@COL"
$wsFixture.Range("L7").Value = "@BOL"
$wsFixture.Range("M7").Value = "VEN"
$wsFixture.Range("N7").Value = "ECU"
$wsFixture.Range("O7").Value = "@BRA"
$wsFixture.Range("P7").Value = "CHI"
$wsFixture.Range("S7").Value = "PAR"
$wsFixture.Range("B8").Value = "@COL"
$wsFixture.Range("C8").Value = "URU"
$wsFixture.Range("D8").Value = "@BRA"
$wsFixture.Range("E8").Value = "PER"
$wsFixture.Range("F8").Value = "@CHI"
$wsFixture.Range("H8").Value = "@BOL"
$wsFixture.Range("I8").Value = "ECU"
$wsFixture.Range("J8").Value = "PAR"
$wsFixture.Range("K8").Value = "@URU"
$wsFixture.Range("L8").Value = "BRA"
$wsFixture.Range("M8").Value = "@PER"
$wsFixture.Range("N8").Value = "CHI"
$wsFixture.Range("P8").Value = "BOL"
$wsFixture.Range("Q8").Value = "@ECU"
$wsFixture.Range("R8").Value = "@PAR"
$wsFixture.Range("S8").Value = "COL"
$wsFixture.Range("B9").Value = "BRA"
$wsFixture.Range("C9").Value = "@ECU"
$wsFixture.Range("D9").Value = "@PER"
$wsFixture.Range("E9").Value = "PAR"
$wsFixture.Range("G9").Value = "URU"
$wsFixture.Range("H9").Value = "VEN"
$wsFixture.Range("J9").Value = "@CHI"
$wsFixture.Range("K9").Value = "ECU"
$wsFixture.Range("L9").Value = "PER"
$wsFixture.Range("M9").Value = "@PAR"
$wsFixture.Range("O9").Value = "@URU"
$wsFixture.Range("P9").Value = "@VEN"
$wsFixture.Range("R9").Value = "CHI"
$wsFixture.Range("S9").Value = "@BRA"
$wsFixture.Range("B10").Value = "PER"
$wsFixture.Range("D10").Value = "COL"
$wsFixture.Range("E10").Value = "@BOL"
$wsFixture.Range("F10").Value = "@BRA"
$wsFixture.Range("G10").Value = "ECU"
$wsFixture.Range("H10").Value = "@URU"
$wsFixture.Range("J10").Value = "@VEN"
$wsFixture.Range("L10").Value = "@COL"
$wsFixture.Range("M10").Value = "BOL"
$wsFixture.Range("N10").Value = "BRA"
$wsFixture.Range("O10").Value = "@ECU"
$wsFixture.Range("P10").Value = "URU"
$wsFixture.Range("R10").Value = "VEN"
$wsFixture.Range("S10").Value = "@PER"
$wsFixture.Range("C11").Value = "BOL"
$wsFixture.Range("D11").Value = "@CHI"
$wsFixture.Range("E11").Value = "URU"
$wsFixture.Range("F11").Value = "PER"
$wsFixture.Range("G11").Value = "@PAR"
$wsFixture.Range("H11").Value = "BRA"
$wsFixture.Range("I11").Value = "@VEN"
$wsFixture.Range("J11").Value = "COL"
$wsFixture.Range("K11").Value = "@BOL"
$wsFixture.Range("L11").Value = "CHI"
$wsFixture.Range("M11").Value = "@URU"
$wsFixture.Range("N11").Value = "@PER"
$wsFixture.Range("O11").Value = "PAR"
$wsFixture.Range("P11").Value = "@BRA"
$wsFixture.Range("Q11").Value = "VEN"
$wsFixture.Range("R11").Value = "@COL"

# --- Sheet "Breaks y secuencias": swap H-A / A-H counts for rows 5 and 8 ---
$wsBreaks = $wb.Worksheets.Item("Breaks y secuencias")
$wsBreaks.Range("C5").Value = 5
$wsBreaks.Range("D5").Value = 4
$wsBreaks.Range("C8").Value = 4
$wsBreaks.Range("D8").Value = 5

# --- Sheet "Partidos acumulados": updated cumulative match counts ---
$wsAcum = $wb.Worksheets.Item("Partidos acumulados")
$wsAcum.Range("F2").Value = 2
$wsAcum.Range("N2").Value = 7
$wsAcum.Range("D4").Value = 2
$wsAcum.Range("L4").Value = 5
$wsAcum.Range("B5").Value = 0
$wsAcum.Range("D5").Value = 1
$wsAcum.Range("J5").Value = 5
$wsAcum.Range("L5").Value = 6
$wsAcum.Range("R5").Value = 8
$wsAcum.Range("D6").Value = 1
$wsAcum.Range("F6").Value = 2
$wsAcum.Range("H6").Value = 3
$wsAcum.Range("L6").Value = 6
$wsAcum.Range("N6").Value = 7
$wsAcum.Range("P6").Value = 8
$wsAcum.Range("F7").Value = 3
$wsAcum.Range("N7").Value = 6
$wsAcum.Range("B8").Value = 1
$wsAcum.Range("H8").Value = 4
$wsAcum.Range("J8").Value = 4
$wsAcum.Range("P8").Value = 7
$wsAcum.Range("R8").Value = 9
$wsAcum.Range("D9").Value = 2
$wsAcum.Range("L9").Value = 5
$wsAcum.Range("D10").Value = 1
$wsAcum.Range("F10").Value = 3
$wsAcum.Range("L10").Value = 6
$wsAcum.Range("N10").Value = 6
$wsAcum.Range("D11").Value = 2
$wsAcum.Range("L11").Value = 5
